$d = $word.ActiveDocument

# Capture the existing (first) paragraph's opening tag verbatim so we keep
# its w14:paraId / rsid* attributes untouched.
$para = $d.Paragraphs(1)
$paraOxml = $para.Range.WordOpenXML
$openTag = '<w:p>'
if ($paraOxml -match '(<w:p[ >][^>]*>)') {
    $openTag = $matches[1]
}

# New run/proofErr content for the paragraph: the plain placeholder text
# "2<10" and "2{{fields.b}}" gets surrounded by the spell-checker's
# spellStart/spellEnd markers (and extra empty runs), and the
# "{{fields.b}}" literal is split into "{{", "fields.b" (flagged), "}}" so
# that "fields.b" renders/parses as its own token.
$body = '<w:r/>' + `
    '<w:proofErr w:type="spellStart"/><w:proofErr w:type="spellEnd"/>' + `
    '<w:proofErr w:type="spellStart"/><w:proofErr w:type="spellEnd"/>' + `
    '<w:proofErr w:type="spellStart"/><w:proofErr w:type="spellEnd"/>' + `
    '<w:proofErr w:type="spellStart"/><w:proofErr w:type="spellEnd"/>' + `
    '<w:proofErr w:type="spellStart"/><w:proofErr w:type="spellEnd"/>' + `
    '<w:proofErr w:type="spellStart"/><w:proofErr w:type="spellEnd"/>' + `
    '<w:r/>' + `
    '<w:proofErr w:type="spellStart"/><w:r><w:t/></w:r><w:proofErr w:type="spellEnd"/>' + `
    '<w:r/>' + `
    '<w:r><w:t/></w:r>' + `
    '<w:proofErr w:type="spellStart"/><w:r><w:t/></w:r><w:proofErr w:type="spellEnd"/>' + `
    '<w:r/>' + `
    '<w:r><w:t/></w:r>' + `
    '<w:proofErr w:type="spellStart"/><w:r><w:t/></w:r><w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t/></w:r>' + `
    '<w:r><w:t>2</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve">&lt;10</w:t></w:r>' + `
    '<w:r/>' + `
    '<w:proofErr w:type="spellStart"/><w:r><w:t/></w:r><w:proofErr w:type="spellEnd"/>' + `
    '<w:r/>' + `
    '<w:r><w:t/></w:r>' + `
    '<w:proofErr w:type="spellStart"/><w:r><w:t/></w:r><w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t/></w:r>' + `
    '<w:r><w:t>2</w:t></w:r>' + `
    '<w:r><w:t>{{</w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/><w:r><w:t>fields.b</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t xml:space="preserve">}}</w:t></w:r>' + `
    '<w:r/>' + `
    '<w:proofErr w:type="spellStart"/><w:r><w:t/></w:r><w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t/></w:r>' + `
    '<w:r/>' + `
    '<w:proofErr w:type="spellStart"/><w:r><w:t/></w:r><w:proofErr w:type="spellEnd"/>' + `
    '<w:r/>' + `
    '<w:proofErr w:type="spellStart"/><w:proofErr w:type="spellEnd"/>' + `
    '<w:r/>' + `
    '<w:proofErr w:type="spellStart"/><w:r><w:t/></w:r><w:proofErr w:type="spellEnd"/>' + `
    '<w:r/>' + `
    '<w:proofErr w:type="spellStart"/><w:proofErr w:type="spellEnd"/>' + `
    '<w:r/>' + `
    '<w:proofErr w:type="spellStart"/><w:r><w:t/></w:r><w:proofErr w:type="spellEnd"/>' + `
    '<w:r/>' + `
    '<w:proofErr w:type="spellStart"/><w:proofErr w:type="spellEnd"/>' + `
    '<w:r/>' + `
    '<w:proofErr w:type="spellStart"/><w:r><w:t/></w:r><w:proofErr w:type="spellEnd"/>' + `
    '<w:r/>' + `
    '<w:proofErr w:type="spellStart"/><w:proofErr w:type="spellEnd"/>' + `
    '<w:r/>' + `
    '<w:proofErr w:type="spellStart"/><w:r><w:t/></w:r><w:proofErr w:type="spellEnd"/>' + `
    '<w:r/>' + `
    '<w:proofErr w:type="spellStart"/><w:proofErr w:type="spellEnd"/>' + `
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' + `
    '<w:r/>' + `
    '<w:proofErr w:type="spellStart"/><w:r><w:t/></w:r><w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t/></w:r>'

$newParaXml = $openTag + $body + '</w:p>'

$pkg = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?>' + `
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
    '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" mc:Ignorable="w14">' + `
    '<w:body>' + $newParaXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# Replace the whole body (paragraph + its end-of-story mark) in one shot so
# we don't leave a stray empty paragraph behind.
$d.Content.InsertXML($pkg)
